$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.096.60"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.818.42"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "450.98"
$ws.Range("E5").Value = "  +7.32%  "
$ws.Range("D6").Value = "147.06"
$ws.Range("E6").Value = "  +14.39%  "
$ws.Range("E7").Value = "  +3.50%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +3.18%  "
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").Value = "0.0000316"
$ws.Range("E11").Value = "  -8.21%  "
$ws.Range("D12").Value = "43.68"
$ws.Range("E12").Value = "  +9.40%  "
$ws.Range("E13").Value = "  +3.07%  "
$ws.Range("D14").Value = "4.433.06"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "14.83"
$ws.Range("E15").Value = "  -5.40%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.860.47"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.137"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("E19").Value = "  +9.28%  "
$ws.Range("D20").Value = "67.187.03"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "429.26"
$ws.Range("E21").Value = "  +6.68%  "
$ws.Range("D22").Value = "14.74"
$ws.Range("E22").Value = "  +4.70%  "
$ws.Range("E23").Value = "  +8.97%  "
$ws.Range("D24").Value = "86.15"
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "3.48"
$ws.Range("E25").Value = "  +9.22%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "10.26"
$ws.Range("E26").Value = "  +22.04%  "
$ws.Range("D27").Value = "37.06"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "9.70"
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("D30").Value = "736.63"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "13.79"
$ws.Range("E31").Value = "  +12.20%  "
$ws.Range("E32").Value = "  +12.36%  "
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").Value = "42.93"
$ws.Range("E34").Value = "  +12.20%  "
$ws.Range("E35").Value = "  +5.34%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "57.21"
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "5.62"
$ws.Range("E37").Value = "  +10.17%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +6.46%  "
$ws.Range("D40").Value = "0.350"
$ws.Range("E40").Value = "  +14.79%  "
$ws.Range("D41").Value = "2.66"
$ws.Range("E41").Value = "  +20.66%  "
$ws.Range("D42").Value = "2.90"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "0.0₃0680"
$ws.Range("E43").Value = "  -9.74%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D46").Value = "3.46"
$ws.Range("E46").Value = "  +4.93%  "
$ws.Range("E47").Value = "  +6.30%  "
$ws.Range("E48").Value = "  +5.61%  "
$ws.Range("D49").Value = "2.67"
$ws.Range("E49").Value = "  +5.07%  "
$ws.Range("D50").Value = "143.59"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("E51").Value = "  +3.66%  "
